# Update "想去人数" (interested count) values in the F column
# on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 242
$ws1.Range("F5").Value = 2852
$ws1.Range("F10").Value = 196
$ws1.Range("F11").Value = 292
$ws1.Range("F12").Value = 44

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 242
$ws4.Range("F5").Value = 2852
$ws4.Range("F11").Value = 196
$ws4.Range("F12").Value = 292
$ws4.Range("F13").Value = 44
